$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 466.5625
$ws.Range("I28").Value = 431
$ws.Range("J28").Value = 1000
$ws.Range("K28").Value = 431
$ws.Range("L28").Value = 1000
$ws.Range("M28").Value = 54
$ws.Range("N28").Value = -1970

$ws.Range("H40").Value = 2016.1666
$ws.Range("J40").Value = 2833.3333
$ws.Range("L40").Value = 2833.3333
$ws.Range("N40").Value = -3183.3333

$ws.Range("H53").Value = 278.31033
$ws.Range("I53").Value = 257.1111
$ws.Range("J53").Value = 313
$ws.Range("K53").Value = 257.1111
$ws.Range("L53").Value = 313
$ws.Range("M53").Value = 379.8889
$ws.Range("N53").Value = -1587

$ws.Range("H107").Value = 1025.7142
$ws.Range("I107").Value = 1104.1428
$ws.Range("J107").Value = 868.8570999999999
$ws.Range("K107").Value = 1104.1428
$ws.Range("L107").Value = 868.8570999999999
$ws.Range("M107").Value = 815.8571999999999
$ws.Range("N107").Value = -4708.8571

$ws.Range("H111").Value = 795.2308
$ws.Range("I111").Value = 767.0909
$ws.Range("K111").Value = 2301.2727
$ws.Range("M111").Value = 765.7273

$ws.Range("H132").Value = 5003022.5
$ws.Range("I132").Value = 5408335.5
$ws.Range("J132").Value = 4166.3335
$ws.Range("K132").Value = 16225006.5
$ws.Range("L132").Value = 12499.0005
$ws.Range("M132").Value = -16222476.5
$ws.Range("N132").Value = -17559.0005

$ws.Range("H137").Value = 3661.1628
$ws.Range("I137").Value = 4198.7036
$ws.Range("J137").Value = 2754.0625
$ws.Range("K137").Value = 12596.1108
$ws.Range("L137").Value = 8262.1875
$ws.Range("M137").Value = -10046.1108
$ws.Range("N137").Value = -13362.1875

$ws.Range("H138").Value = 5038.8213
$ws.Range("J138").Value = 5832.453
$ws.Range("L138").Value = 17497.359
$ws.Range("N138").Value = -27777.359

$ws.Range("H141").Value = 274296.84
$ws.Range("I141").Value = 955.2
$ws.Range("J141").Value = 537125.4
$ws.Range("K141").Value = 2865.6
$ws.Range("L141").Value = 1611376.2
$ws.Range("M141").Value = 2314.4
$ws.Range("N141").Value = -1621736.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 10000
$ws.Range("I8").Value = 8000
$ws.Range("J8").Value = 11000
$ws.Range("K8").Value = 8000
$ws.Range("L8").Value = 11000
$ws.Range("M8").Value = -7856
$ws.Range("N8").Value = -11288

$ws.Range("H13").Value = 100004
$ws.Range("J13").Value = 100004
$ws.Range("L13").Value = 100004
$ws.Range("N13").Value = -100292

$ws.Range("H32").Value = 6515.61
$ws.Range("I32").Value = 6515.61
$ws.Range("K32").Value = 6515.61
$ws.Range("M32").Value = -6228.61

$ws.Range("H74").Value = 1792.4878
$ws.Range("I74").Value = 1412.4117
$ws.Range("K74").Value = 1412.4117
$ws.Range("M74").Value = -538.4117000000001

$ws.Range("H77").Value = 1792.4878
$ws.Range("I77").Value = 1412.4117
$ws.Range("K77").Value = 7062.058500000001
$ws.Range("M77").Value = -2694.058500000001

$ws.Range("H97").Value = 900
$ws.Range("I97").Value = 750
$ws.Range("J97").Value = 1200
$ws.Range("K97").Value = 750
$ws.Range("L97").Value = 1200
$ws.Range("M97").Value = -254
$ws.Range("N97").Value = -2192

$ws.Range("H102").Value = 3659.9
$ws.Range("I102").Value = 1646.9231
$ws.Range("J102").Value = 7398.2856
$ws.Range("K102").Value = 1646.9231
$ws.Range("L102").Value = 7398.2856
$ws.Range("M102").Value = -24.92309999999998
$ws.Range("N102").Value = -10642.2856

$ws.Range("H110").Value = 1591.4828
$ws.Range("I110").Value = 639.2083
$ws.Range("J110").Value = 6162.4
$ws.Range("K110").Value = 639.2083
$ws.Range("L110").Value = 6162.4
$ws.Range("M110").Value = 1405.7917
$ws.Range("N110").Value = -10252.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 855.5
$ws.Range("I94").Value = 612.8
$ws.Range("J94").Value = 1664.5
$ws.Range("K94").Value = 612.8
$ws.Range("L94").Value = 1664.5
$ws.Range("M94").Value = -161.8
$ws.Range("N94").Value = -2566.5

$ws.Range("H107").Value = 1704.1852
$ws.Range("I107").Value = 1424.25
$ws.Range("J107").Value = 2111.3635
$ws.Range("K107").Value = 1424.25
$ws.Range("L107").Value = 2111.3635
$ws.Range("M107").Value = 495.75
$ws.Range("N107").Value = -5951.363499999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 8067521.5
$ws.Range("I58").Value = 1622.2162
$ws.Range("J58").Value = 20005052
$ws.Range("K58").Value = 1622.2162
$ws.Range("L58").Value = 20005052
$ws.Range("M58").Value = -1419.2162
$ws.Range("N58").Value = -20005458

$ws.Range("H136").Value = 8067521.5
$ws.Range("I136").Value = 1622.2162
$ws.Range("J136").Value = 20005052
$ws.Range("K136").Value = 4866.6486
$ws.Range("L136").Value = 60015156
$ws.Range("M136").Value = -2316.6486
$ws.Range("N136").Value = -60020256

$ws.Range("H141").Value = 25094.445
$ws.Range("J141").Value = 25094.445
$ws.Range("L141").Value = 25094.445
$ws.Range("N141").Value = -35454.445

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H120").Value = 17211.428
$ws.Range("I120").Value = 3740
$ws.Range("K120").Value = 11220
$ws.Range("M120").Value = -6382

$ws.Range("H131").Value = 1253.7013
$ws.Range("J131").Value = 1033.2239
$ws.Range("L131").Value = 3099.6717
$ws.Range("N131").Value = -13179.6717

$ws.Range("H132").Value = 2802.5
$ws.Range("I132").Value = 1828.5714
$ws.Range("J132").Value = 3560
$ws.Range("K132").Value = 16457.1426
$ws.Range("L132").Value = 32040
$ws.Range("M132").Value = -13927.1426
$ws.Range("N132").Value = -37100

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2129.125
$ws.Range("I113").Value = 1004.7143
$ws.Range("K113").Value = 1004.7143
$ws.Range("M113").Value = 1165.2857

$ws.Range("H132").Value = 4804.343
$ws.Range("I132").Value = 5477.7
$ws.Range("J132").Value = 3906.5334
$ws.Range("K132").Value = 16433.1
$ws.Range("L132").Value = 11719.6002
$ws.Range("M132").Value = -13903.1
$ws.Range("N132").Value = -16779.6002

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 142864540
$ws.Range("I61").Value = 333336670
$ws.Range("J61").Value = 10450
$ws.Range("K61").Value = 333336670
$ws.Range("L61").Value = 10450
$ws.Range("M61").Value = -333336468
$ws.Range("N61").Value = -10854

$ws.Range("H93").Value = 2799.875
$ws.Range("I93").Value = 1816.5
$ws.Range("J93").Value = 5750
$ws.Range("K93").Value = 1816.5
$ws.Range("L93").Value = 5750
$ws.Range("M93").Value = -568.5
$ws.Range("N93").Value = -8246

$ws.Range("H113").Value = 142864540
$ws.Range("I113").Value = 333336670
$ws.Range("J113").Value = 10450
$ws.Range("K113").Value = 333336670
$ws.Range("L113").Value = 10450
$ws.Range("M113").Value = -333334500
$ws.Range("N113").Value = -14790

$ws.Range("H132").Value = 2492.6
$ws.Range("I132").Value = 1763.1072
$ws.Range("J132").Value = 3694.1177
$ws.Range("K132").Value = 5289.321599999999
$ws.Range("L132").Value = 11082.3531
$ws.Range("M132").Value = -2759.321599999999
$ws.Range("N132").Value = -16142.3531

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 29933.334
$ws.Range("J92").Value = 29933.334
$ws.Range("L92").Value = 29933.334
$ws.Range("N92").Value = -34925.334
